$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.079437672839434
$ws.Cells.Item(2, 4).Value = 1.082468570995508
$ws.Cells.Item(2, 5).Value = 1.077609996695981
$ws.Cells.Item(2, 6).Value = 1.091387622825828
$ws.Cells.Item(2, 9).Value = 1.064830520483784
$ws.Cells.Item(2, 10).Value = 1.084322175320534
$ws.Cells.Item(2, 11).Value = 1.0851367087605
$ws.Cells.Item(2, 12).Value = 1.080290844615396
$ws.Cells.Item(2, 13).Value = 1.094032745868362
$ws.Cells.Item(2, 14).Value = 1.085862037182981

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.080818345942717
$ws.Cells.Item(3, 4).Value = 1.083622919188407
$ws.Cells.Item(3, 5).Value = 1.078919175726516
$ws.Cells.Item(3, 6).Value = 1.092725114617711
$ws.Cells.Item(3, 9).Value = 1.065359735033636
$ws.Cells.Item(3, 10).Value = 1.085361809124349
$ws.Cells.Item(3, 11).Value = 1.08610888314138
$ws.Cells.Item(3, 12).Value = 1.081416550490377
$ws.Cells.Item(3, 13).Value = 1.095189297141629
$ws.Cells.Item(3, 14).Value = 1.086903147386045

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.08171040081098
$ws.Cells.Item(4, 4).Value = 1.084368551582291
$ws.Cells.Item(4, 5).Value = 1.079763922590518
$ws.Cells.Item(4, 6).Value = 1.093588958903486
$ws.Cells.Item(4, 9).Value = 1.065700158957378
$ws.Cells.Item(4, 10).Value = 1.086032677944361
$ws.Cells.Item(4, 11).Value = 1.08673605182142
$ws.Cells.Item(4, 12).Value = 1.082142028096045
$ws.Cells.Item(4, 13).Value = 1.095935510362141
$ws.Cells.Item(4, 14).Value = 1.087574968916731

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.082085107813336
$ws.Cells.Item(5, 4).Value = 1.084681707093636
$ws.Cells.Item(5, 5).Value = 1.080118491681366
$ws.Cells.Item(5, 6).Value = 1.093951741803806
$ws.Cells.Item(5, 9).Value = 1.065842794029968
$ws.Cells.Item(5, 10).Value = 1.086314274674831
$ws.Cells.Item(5, 11).Value = 1.086999264907732
$ws.Cells.Item(5, 12).Value = 1.08244632426267
$ws.Cells.Item(5, 13).Value = 1.09624870896518
$ws.Cells.Item(5, 14).Value = 1.087856965546858

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.082148004602233
$ws.Cells.Item(6, 4).Value = 1.084734269339517
$ws.Cells.Item(6, 5).Value = 1.080177992614787
$ws.Cells.Item(6, 6).Value = 1.094012632712048
$ws.Cells.Item(6, 9).Value = 1.06586671507113
$ws.Cells.Item(6, 10).Value = 1.086361530548425
$ws.Cells.Item(6, 11).Value = 1.087043433368514
$ws.Cells.Item(6, 12).Value = 1.082497376388636
$ws.Cells.Item(6, 13).Value = 1.096301266687483
$ws.Cells.Item(6, 14).Value = 1.087904288529213

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.08171540888588
$ws.Cells.Item(7, 4).Value = 1.084372737187198
$ws.Cells.Item(7, 5).Value = 1.079768662561373
$ws.Cells.Item(7, 6).Value = 1.093593807900745
$ws.Cells.Item(7, 9).Value = 1.065702066733405
$ws.Cells.Item(7, 10).Value = 1.086036442362846
$ws.Cells.Item(7, 11).Value = 1.086739570644669
$ws.Cells.Item(7, 12).Value = 1.082146096837416
$ws.Cells.Item(7, 13).Value = 1.095939697331564
$ws.Cells.Item(7, 14).Value = 1.087578738681122

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.079904557545791
$ws.Cells.Item(8, 4).Value = 1.082858962121182
$ws.Cells.Item(8, 5).Value = 1.0780529368987
$ws.Cells.Item(8, 6).Value = 1.091839970055864
$ws.Cells.Item(8, 9).Value = 1.065009789594304
$ws.Cells.Item(8, 10).Value = 1.084673909142955
$ws.Cells.Item(8, 11).Value = 1.085465654409366
$ws.Cells.Item(8, 12).Value = 1.080671892761621
$ws.Cells.Item(8, 13).Value = 1.094424057124838
$ws.Cells.Item(8, 14).Value = 1.086214270507774

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.076703116447553
$ws.Cells.Item(9, 4).Value = 1.080181262262247
$ws.Cells.Item(9, 5).Value = 1.075011062679034
$ws.Cells.Item(9, 6).Value = 1.088736925488656
$ws.Cells.Item(9, 9).Value = 1.063774361789257
$ws.Cells.Item(9, 10).Value = 1.082258610862702
$ws.Cells.Item(9, 11).Value = 1.083206145064131
$ws.Cells.Item(9, 12).Value = 1.078051421255545
$ws.Cells.Item(9, 13).Value = 1.091736550251958
$ws.Cells.Item(9, 14).Value = 1.083795542226925

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.074561374845161
$ws.Cells.Item(10, 4).Value = 1.078388942627017
$ws.Cells.Item(10, 5).Value = 1.072970224560134
$ws.Cells.Item(10, 6).Value = 1.086659383340199
$ws.Cells.Item(10, 9).Value = 1.06294011029789
$ws.Cells.Item(10, 10).Value = 1.080638461014891
$ws.Cells.Item(10, 11).Value = 1.081689626545275
$ws.Cells.Item(10, 12).Value = 1.07628873275704
$ws.Cells.Item(10, 13).Value = 1.089933237912336
$ws.Cells.Item(10, 14).Value = 1.082173091580496

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.073632122113622
$ws.Cells.Item(11, 4).Value = 1.077611075305607
$ws.Cells.Item(11, 5).Value = 1.072083352953826
$ws.Cells.Item(11, 6).Value = 1.085757596381488
$ws.Cells.Item(11, 9).Value = 1.062576307708906
$ws.Cells.Item(11, 10).Value = 1.079934488576409
$ws.Cells.Item(11, 11).Value = 1.081030477515419
$ws.Cells.Item(11, 12).Value = 1.075521649273864
$ws.Cells.Item(11, 13).Value = 1.089149538418754
$ws.Cells.Item(11, 14).Value = 1.081468119420405

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.073286668444973
$ws.Cells.Item(12, 4).Value = 1.077321867157773
$ws.Cells.Item(12, 5).Value = 1.071753443497181
$ws.Cells.Item(12, 6).Value = 1.085422294548085
$ws.Cells.Item(12, 9).Value = 1.062440786132543
$ws.Cells.Item(12, 10).Value = 1.079672629689012
$ws.Cells.Item(12, 11).Value = 1.080785260922612
$ws.Cells.Item(12, 12).Value = 1.075236137027299
$ws.Cells.Item(12, 13).Value = 1.088858001504812
$ws.Cells.Item(12, 14).Value = 1.081205888663352

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.073360782578882
$ws.Cells.Item(13, 4).Value = 1.077383915764636
$ws.Cells.Item(13, 5).Value = 1.071824232323633
$ws.Cells.Item(13, 6).Value = 1.085494233351763
$ws.Cells.Item(13, 9).Value = 1.062469873635535
$ws.Cells.Item(13, 10).Value = 1.079728816249101
$ws.Cells.Item(13, 11).Value = 1.080837877970464
$ws.Cells.Item(13, 12).Value = 1.075297406900838
$ws.Cells.Item(13, 13).Value = 1.088920557012185
$ws.Cells.Item(13, 14).Value = 1.081262155014801

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.073603572722048
$ws.Cells.Item(14, 4).Value = 1.077587174885934
$ws.Cells.Item(14, 5).Value = 1.0720560924977
$ws.Cells.Item(14, 6).Value = 1.085729887168192
$ws.Cells.Item(14, 9).Value = 1.062565113411246
$ws.Cells.Item(14, 10).Value = 1.079912850866708
$ws.Cells.Item(14, 11).Value = 1.081010215606136
$ws.Cells.Item(14, 12).Value = 1.075498060689086
$ws.Cells.Item(14, 13).Value = 1.089125448824016
$ws.Cells.Item(14, 14).Value = 1.081446450982674

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.073753125370497
$ws.Cells.Item(15, 4).Value = 1.077712373105245
$ws.Cells.Item(15, 5).Value = 1.072198884558206
$ws.Cells.Item(15, 6).Value = 1.0858750362371
$ws.Cells.Item(15, 9).Value = 1.062623742114461
$ws.Cells.Item(15, 10).Value = 1.080026191168795
$ws.Cells.Item(15, 11).Value = 1.081116348079073
$ws.Cells.Item(15, 12).Value = 1.075621612572593
$ws.Cells.Item(15, 13).Value = 1.089251631433004
$ws.Cells.Item(15, 14).Value = 1.081559952240991

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.074623006349098
$ws.Cells.Item(16, 4).Value = 1.078440529082628
$ws.Cells.Item(16, 5).Value = 1.073029015658696
$ws.Cells.Item(16, 6).Value = 1.086719185002768
$ws.Cells.Item(16, 9).Value = 1.062964200308064
$ws.Cells.Item(16, 10).Value = 1.080685129488383
$ws.Cells.Item(16, 11).Value = 1.081733319220098
$ws.Cells.Item(16, 12).Value = 1.076339560215917
$ws.Cells.Item(16, 13).Value = 1.089985188724258
$ws.Cells.Item(16, 14).Value = 1.082219826328573

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.075168155103249
$ws.Cells.Item(17, 4).Value = 1.078896801019117
$ws.Cells.Item(17, 5).Value = 1.073548878665773
$ws.Cells.Item(17, 6).Value = 1.087248103361106
$ws.Cells.Item(17, 9).Value = 1.063177071194934
$ws.Cells.Item(17, 10).Value = 1.081097807400472
$ws.Cells.Item(17, 11).Value = 1.082119659180962
$ws.Cells.Item(17, 12).Value = 1.076788879232399
$ws.Cells.Item(17, 13).Value = 1.09044456074198
$ws.Cells.Item(17, 14).Value = 1.082633090290632

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.075485951777158
$ws.Cells.Item(18, 4).Value = 1.079162765279213
$ws.Cells.Item(18, 5).Value = 1.073851800197434
$ws.Cells.Item(18, 6).Value = 1.087556400998461
$ws.Cells.Item(18, 9).Value = 1.063300987823803
$ws.Cells.Item(18, 10).Value = 1.081338280735613
$ws.Cells.Item(18, 11).Value = 1.082344765221866
$ws.Cells.Item(18, 12).Value = 1.077050590650517
$ws.Cells.Item(18, 13).Value = 1.090712229860451
$ws.Cells.Item(18, 14).Value = 1.082873905125492

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.075594282038029
$ws.Cells.Item(19, 4).Value = 1.079253423319937
$ws.Cells.Item(19, 5).Value = 1.073955037013849
$ws.Cells.Item(19, 6).Value = 1.087661486889366
$ws.Cells.Item(19, 9).Value = 1.06334319834275
$ws.Cells.Item(19, 10).Value = 1.081420236300399
$ws.Cells.Item(19, 11).Value = 1.082421480109765
$ws.Cells.Item(19, 12).Value = 1.077139765179752
$ws.Cells.Item(19, 13).Value = 1.09080345172165
$ws.Cells.Item(19, 14).Value = 1.082955977076581

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.07510968441188
$ws.Cells.Item(20, 4).Value = 1.078847865136225
$ws.Cells.Item(20, 5).Value = 1.073493133977309
$ws.Cells.Item(20, 6).Value = 1.087191377361089
$ws.Cells.Item(20, 9).Value = 1.063154257775317
$ws.Cells.Item(20, 10).Value = 1.081053555288396
$ws.Cells.Item(20, 11).Value = 1.08207823337706
$ws.Cells.Item(20, 12).Value = 1.076740709758162
$ws.Cells.Item(20, 13).Value = 1.090395302951281
$ws.Cells.Item(20, 14).Value = 1.082588775335482

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.073532085083986
$ws.Cells.Item(21, 4).Value = 1.077527327741681
$ws.Cells.Item(21, 5).Value = 1.071987828938441
$ws.Cells.Item(21, 6).Value = 1.085660502386746
$ws.Cells.Item(21, 9).Value = 1.062537078440751
$ws.Cells.Item(21, 10).Value = 1.079858667581709
$ws.Cells.Item(21, 11).Value = 1.080959476984426
$ws.Cells.Item(21, 12).Value = 1.075438989317433
$ws.Cells.Item(21, 13).Value = 1.089065125387985
$ws.Cells.Item(21, 14).Value = 1.081392190751195

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.072538516620261
$ws.Cells.Item(22, 4).Value = 1.076695467554206
$ws.Cells.Item(22, 5).Value = 1.071038568502557
$ws.Cells.Item(22, 6).Value = 1.084696021004768
$ws.Cells.Item(22, 9).Value = 1.062146780324079
$ws.Cells.Item(22, 10).Value = 1.079105237901011
$ws.Cells.Item(22, 11).Value = 1.080253873087883
$ws.Cells.Item(22, 12).Value = 1.074617167184938
$ws.Cells.Item(22, 13).Value = 1.088226263560108
$ws.Cells.Item(22, 14).Value = 1.080637691113932

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.073065386867084
$ws.Cells.Item(23, 4).Value = 1.077136604721187
$ws.Cells.Item(23, 5).Value = 1.071542059246873
$ws.Cells.Item(23, 6).Value = 1.085207499211902
$ws.Cells.Item(23, 9).Value = 1.06235389949026
$ws.Cells.Item(23, 10).Value = 1.079504851665413
$ws.Cells.Item(23, 11).Value = 1.080628137270031
$ws.Cells.Item(23, 12).Value = 1.07505315367746
$ws.Cells.Item(23, 13).Value = 1.088671202026037
$ws.Cells.Item(23, 14).Value = 1.081037872375715

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.075136105354105
$ws.Cells.Item(24, 4).Value = 1.078869977684572
$ws.Cells.Item(24, 5).Value = 1.073518323546568
$ws.Cells.Item(24, 6).Value = 1.087217010051621
$ws.Cells.Item(24, 9).Value = 1.063164566941668
$ws.Cells.Item(24, 10).Value = 1.08107355163778
$ws.Cells.Item(24, 11).Value = 1.082096952654335
$ws.Cells.Item(24, 12).Value = 1.076762476608869
$ws.Cells.Item(24, 13).Value = 1.090417561275145
$ws.Cells.Item(24, 14).Value = 1.082608800081975

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.077532049271891
$ws.Cells.Item(25, 4).Value = 1.080874755997059
$ws.Cells.Item(25, 5).Value = 1.075799704306252
$ws.Cells.Item(25, 6).Value = 1.089540666345601
$ws.Cells.Item(25, 9).Value = 1.064095610341376
$ws.Cells.Item(25, 10).Value = 1.082884754488533
$ws.Cells.Item(25, 11).Value = 1.083792055005916
$ws.Cells.Item(25, 12).Value = 1.078731613439446
$ws.Cells.Item(25, 13).Value = 1.092433359303651
$ws.Cells.Item(25, 14).Value = 1.084422575048525
